$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (ECs -> Ntn4 -> Dcc -> FAPs) ---
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 2.250039333333334
$ws.Range("H2").Value = 6.750118000000001
$ws.Range("I2").Value = 0.05266340474459803
$ws.Range("J2").Value = 0.05266340474459803
$ws.Range("M2").Value = 0.05453333333333333
$ws.Range("N2").Value = 0.1636
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.1227021449777778
$ws.Range("R2").Value = 1.1043193048
$ws.Range("S2").Value = 0.05266340474459803
$ws.Range("T2").Value = 0.05266340474459803

# --- Row 3 (FAPs -> Ntn4 -> Dcc -> FAPs) ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 18.996877
$ws.Range("H3").Value = 56.99063100000001
$ws.Range("I3").Value = 0.4446323259834918
$ws.Range("J3").Value = 0.4446323259834918
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1.035963025733333
$ws.Range("R3").Value = 9.323667231600002
$ws.Range("S3").Value = 0.4446323259834918
$ws.Range("T3").Value = 0.4446323259834918

# --- Row 4 (MuSCs -> Ntn4 -> Dcc -> FAPs) ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 21.38027566666667
$ws.Range("H4").Value = 64.140827
$ws.Range("I4").Value = 0.5004170790724312
$ws.Range("J4").Value = 0.5004170790724313
$ws.Range("M4").Value = 0.05453333333333333
$ws.Range("N4").Value = 0.1636
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 1.165937699688889
$ws.Range("R4").Value = 10.4934392972
$ws.Range("S4").Value = 0.5004170790724312
$ws.Range("T4").Value = 0.5004170790724313

# --- Row 5 (Resolving-Mac -> Ntn4 -> Dcc -> FAPs) ---
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.09771999999999999
$ws.Range("H5").Value = 0.29316
$ws.Range("I5").Value = 0.002287190199478936
$ws.Range("J5").Value = 0.002287190199478936
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.005328997333333333
$ws.Range("R5").Value = 0.047960976
$ws.Range("S5").Value = 0.002287190199478936
$ws.Range("T5").Value = 0.002287190199478936

# --- Remove old rows 6-9 (MuSCs/Resolving-Mac duplicate pairs, now obsolete) ---
$ws.Rows("6:9").Delete()
